# Adds a new weekly data point (rows 205-206) and shifts the
# existing historical rows down by two positions, which also
# extends the table by one more historical pair at the bottom
# (rows 326-327). Mirrors the upstream daily-logic refresh for
# this sub-sheet (Vega Monumental Concepción - Brócoli).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 123,18
$arr[0,0] = 11
$arr[0,1] = 'Vega Monumental Concepción'
$arr[0,2] = 'Bíobío'
$arr[0,3] = 44777
$arr[0,4] = 8
$arr[0,5] = 100112023
$arr[0,6] = 'Brócoli'
$arr[0,7] = 'Sin especificar'
$arr[0,8] = 'Primera'
$arr[0,9] = 2000
$arr[0,10] = 800
$arr[0,11] = 900
$arr[0,12] = 850
$arr[0,13] = '$/unidad'
$arr[0,14] = 'Región Metropolitana'
$arr[0,15] = 850
$arr[0,16] = 1
$arr[0,17] = 'Hortaliza'
$arr[1,0] = 11
$arr[1,1] = 'Vega Monumental Concepción'
$arr[1,2] = 'Bíobío'
$arr[1,3] = 44777
$arr[1,4] = 8
$arr[1,5] = 100112023
$arr[1,6] = 'Brócoli'
$arr[1,7] = 'Sin especificar'
$arr[1,8] = 'Segunda'
$arr[1,9] = 1500
$arr[1,10] = 600
$arr[1,11] = 600
$arr[1,12] = 600
$arr[1,13] = '$/unidad'
$arr[1,14] = 'Región Metropolitana'
$arr[1,15] = 600
$arr[1,16] = 1
$arr[1,17] = 'Hortaliza'
$arr[2,0] = 11
$arr[2,1] = 'Vega Monumental Concepción'
$arr[2,2] = 'Bíobío'
$arr[2,3] = 44643
$arr[2,4] = 8
$arr[2,5] = 100112023
$arr[2,6] = 'Brócoli'
$arr[2,7] = 'Sin especificar'
$arr[2,8] = 'Primera'
$arr[2,9] = 600
$arr[2,10] = 1000
$arr[2,11] = 1000
$arr[2,12] = 1000
$arr[2,13] = '$/unidad'
$arr[2,14] = 'Región Metropolitana'
$arr[2,15] = 1000
$arr[2,16] = 1
$arr[2,17] = 'Hortaliza'
$arr[3,0] = 11
$arr[3,1] = 'Vega Monumental Concepción'
$arr[3,2] = 'Bíobío'
$arr[3,3] = 44643
$arr[3,4] = 8
$arr[3,5] = 100112023
$arr[3,6] = 'Brócoli'
$arr[3,7] = 'Sin especificar'
$arr[3,8] = 'Segunda'
$arr[3,9] = 500
$arr[3,10] = 800
$arr[3,11] = 800
$arr[3,12] = 800
$arr[3,13] = '$/unidad'
$arr[3,14] = 'Región Metropolitana'
$arr[3,15] = 800
$arr[3,16] = 1
$arr[3,17] = 'Hortaliza'
$arr[4,0] = 11
$arr[4,1] = 'Vega Monumental Concepción'
$arr[4,2] = 'Bíobío'
$arr[4,3] = 44533
$arr[4,4] = 8
$arr[4,5] = 100112023
$arr[4,6] = 'Brócoli'
$arr[4,7] = 'Sin especificar'
$arr[4,8] = 'Primera'
$arr[4,9] = 1500
$arr[4,10] = 600
$arr[4,11] = 650
$arr[4,12] = 623
$arr[4,13] = '$/unidad'
$arr[4,14] = 'Región Metropolitana'
$arr[4,15] = 623
$arr[4,16] = 1
$arr[4,17] = 'Hortaliza'
$arr[5,0] = 11
$arr[5,1] = 'Vega Monumental Concepción'
$arr[5,2] = 'Bíobío'
$arr[5,3] = 44687
$arr[5,4] = 8
$arr[5,5] = 100112023
$arr[5,6] = 'Brócoli'
$arr[5,7] = 'Sin especificar'
$arr[5,8] = 'Primera'
$arr[5,9] = 2000
$arr[5,10] = 700
$arr[5,11] = 800
$arr[5,12] = 750
$arr[5,13] = '$/unidad'
$arr[5,14] = 'Región Metropolitana'
$arr[5,15] = 750
$arr[5,16] = 1
$arr[5,17] = 'Hortaliza'
$arr[6,0] = 11
$arr[6,1] = 'Vega Monumental Concepción'
$arr[6,2] = 'Bíobío'
$arr[6,3] = 44687
$arr[6,4] = 8
$arr[6,5] = 100112023
$arr[6,6] = 'Brócoli'
$arr[6,7] = 'Sin especificar'
$arr[6,8] = 'Segunda'
$arr[6,9] = 2200
$arr[6,10] = 600
$arr[6,11] = 650
$arr[6,12] = 627
$arr[6,13] = '$/unidad'
$arr[6,14] = 'Región Metropolitana'
$arr[6,15] = 627
$arr[6,16] = 1
$arr[6,17] = 'Hortaliza'
$arr[7,0] = 11
$arr[7,1] = 'Vega Monumental Concepción'
$arr[7,2] = 'Bíobío'
$arr[7,3] = 44222
$arr[7,4] = 8
$arr[7,5] = 100112023
$arr[7,6] = 'Brócoli'
$arr[7,7] = 'Sin especificar'
$arr[7,8] = 'Primera'
$arr[7,9] = 1000
$arr[7,10] = 700
$arr[7,11] = 800
$arr[7,12] = 750
$arr[7,13] = '$/unidad'
$arr[7,14] = 'Región del Maule'
$arr[7,15] = 750
$arr[7,16] = 1
$arr[7,17] = 'Hortaliza'
$arr[8,0] = 11
$arr[8,1] = 'Vega Monumental Concepción'
$arr[8,2] = 'Bíobío'
$arr[8,3] = 44222
$arr[8,4] = 8
$arr[8,5] = 100112023
$arr[8,6] = 'Brócoli'
$arr[8,7] = 'Sin especificar'
$arr[8,8] = 'Segunda'
$arr[8,9] = 500
$arr[8,10] = 600
$arr[8,11] = 600
$arr[8,12] = 600
$arr[8,13] = '$/unidad'
$arr[8,14] = 'Región del Maule'
$arr[8,15] = 600
$arr[8,16] = 1
$arr[8,17] = 'Hortaliza'
$arr[9,0] = 11
$arr[9,1] = 'Vega Monumental Concepción'
$arr[9,2] = 'Bíobío'
$arr[9,3] = 44601
$arr[9,4] = 8
$arr[9,5] = 100112023
$arr[9,6] = 'Brócoli'
$arr[9,7] = 'Sin especificar'
$arr[9,8] = 'Primera'
$arr[9,9] = 1500
$arr[9,10] = 750
$arr[9,11] = 800
$arr[9,12] = 773
$arr[9,13] = '$/unidad'
$arr[9,14] = 'Región del Maule'
$arr[9,15] = 773
$arr[9,16] = 1
$arr[9,17] = 'Hortaliza'
$arr[10,0] = 11
$arr[10,1] = 'Vega Monumental Concepción'
$arr[10,2] = 'Bíobío'
$arr[10,3] = 44749
$arr[10,4] = 8
$arr[10,5] = 100112023
$arr[10,6] = 'Brócoli'
$arr[10,7] = 'Sin especificar'
$arr[10,8] = 'Primera'
$arr[10,9] = 3000
$arr[10,10] = 800
$arr[10,11] = 900
$arr[10,12] = 850
$arr[10,13] = '$/unidad'
$arr[10,14] = 'Región Metropolitana'
$arr[10,15] = 850
$arr[10,16] = 1
$arr[10,17] = 'Hortaliza'
$arr[11,0] = 11
$arr[11,1] = 'Vega Monumental Concepción'
$arr[11,2] = 'Bíobío'
$arr[11,3] = 44749
$arr[11,4] = 8
$arr[11,5] = 100112023
$arr[11,6] = 'Brócoli'
$arr[11,7] = 'Sin especificar'
$arr[11,8] = 'Segunda'
$arr[11,9] = 2000
$arr[11,10] = 600
$arr[11,11] = 600
$arr[11,12] = 600
$arr[11,13] = '$/unidad'
$arr[11,14] = 'Región Metropolitana'
$arr[11,15] = 600
$arr[11,16] = 1
$arr[11,17] = 'Hortaliza'
$arr[12,0] = 11
$arr[12,1] = 'Vega Monumental Concepción'
$arr[12,2] = 'Bíobío'
$arr[12,3] = 44161
$arr[12,4] = 8
$arr[12,5] = 100112023
$arr[12,6] = 'Brócoli'
$arr[12,7] = 'Sin especificar'
$arr[12,8] = 'Primera'
$arr[12,9] = 800
$arr[12,10] = 700
$arr[12,11] = 800
$arr[12,12] = 750
$arr[12,13] = '$/unidad'
$arr[12,14] = 'Región del Maule'
$arr[12,15] = 750
$arr[12,16] = 1
$arr[12,17] = 'Hortaliza'
$arr[13,0] = 11
$arr[13,1] = 'Vega Monumental Concepción'
$arr[13,2] = 'Bíobío'
$arr[13,3] = 44161
$arr[13,4] = 8
$arr[13,5] = 100112023
$arr[13,6] = 'Brócoli'
$arr[13,7] = 'Sin especificar'
$arr[13,8] = 'Segunda'
$arr[13,9] = 400
$arr[13,10] = 600
$arr[13,11] = 600
$arr[13,12] = 600
$arr[13,13] = '$/unidad'
$arr[13,14] = 'Región del Maule'
$arr[13,15] = 600
$arr[13,16] = 1
$arr[13,17] = 'Hortaliza'
$arr[14,0] = 11
$arr[14,1] = 'Vega Monumental Concepción'
$arr[14,2] = 'Bíobío'
$arr[14,3] = 44614
$arr[14,4] = 8
$arr[14,5] = 100112023
$arr[14,6] = 'Brócoli'
$arr[14,7] = 'Sin especificar'
$arr[14,8] = 'Primera'
$arr[14,9] = 2000
$arr[14,10] = 900
$arr[14,11] = 1000
$arr[14,12] = 950
$arr[14,13] = '$/unidad'
$arr[14,14] = 'Región Metropolitana'
$arr[14,15] = 950
$arr[14,16] = 1
$arr[14,17] = 'Hortaliza'
$arr[15,0] = 11
$arr[15,1] = 'Vega Monumental Concepción'
$arr[15,2] = 'Bíobío'
$arr[15,3] = 44614
$arr[15,4] = 8
$arr[15,5] = 100112023
$arr[15,6] = 'Brócoli'
$arr[15,7] = 'Sin especificar'
$arr[15,8] = 'Segunda'
$arr[15,9] = 1000
$arr[15,10] = 800
$arr[15,11] = 800
$arr[15,12] = 800
$arr[15,13] = '$/unidad'
$arr[15,14] = 'Región Metropolitana'
$arr[15,15] = 800
$arr[15,16] = 1
$arr[15,17] = 'Hortaliza'
$arr[16,0] = 11
$arr[16,1] = 'Vega Monumental Concepción'
$arr[16,2] = 'Bíobío'
$arr[16,3] = 44420
$arr[16,4] = 8
$arr[16,5] = 100112023
$arr[16,6] = 'Brócoli'
$arr[16,7] = 'Sin especificar'
$arr[16,8] = 'Primera'
$arr[16,9] = 2000
$arr[16,10] = 600
$arr[16,11] = 700
$arr[16,12] = 650
$arr[16,13] = '$/unidad'
$arr[16,14] = 'Región Metropolitana'
$arr[16,15] = 650
$arr[16,16] = 1
$arr[16,17] = 'Hortaliza'
$arr[17,0] = 11
$arr[17,1] = 'Vega Monumental Concepción'
$arr[17,2] = 'Bíobío'
$arr[17,3] = 44420
$arr[17,4] = 8
$arr[17,5] = 100112023
$arr[17,6] = 'Brócoli'
$arr[17,7] = 'Sin especificar'
$arr[17,8] = 'Segunda'
$arr[17,9] = 1000
$arr[17,10] = 500
$arr[17,11] = 500
$arr[17,12] = 500
$arr[17,13] = '$/unidad'
$arr[17,14] = 'Región Metropolitana'
$arr[17,15] = 500
$arr[17,16] = 1
$arr[17,17] = 'Hortaliza'
$arr[18,0] = 11
$arr[18,1] = 'Vega Monumental Concepción'
$arr[18,2] = 'Bíobío'
$arr[18,3] = 44195
$arr[18,4] = 8
$arr[18,5] = 100112023
$arr[18,6] = 'Brócoli'
$arr[18,7] = 'Sin especificar'
$arr[18,8] = 'Primera'
$arr[18,9] = 600
$arr[18,10] = 700
$arr[18,11] = 800
$arr[18,12] = 750
$arr[18,13] = '$/unidad'
$arr[18,14] = 'Región del Maule'
$arr[18,15] = 750
$arr[18,16] = 1
$arr[18,17] = 'Hortaliza'
$arr[19,0] = 11
$arr[19,1] = 'Vega Monumental Concepción'
$arr[19,2] = 'Bíobío'
$arr[19,3] = 44195
$arr[19,4] = 8
$arr[19,5] = 100112023
$arr[19,6] = 'Brócoli'
$arr[19,7] = 'Sin especificar'
$arr[19,8] = 'Segunda'
$arr[19,9] = 300
$arr[19,10] = 600
$arr[19,11] = 600
$arr[19,12] = 600
$arr[19,13] = '$/unidad'
$arr[19,14] = 'Región del Maule'
$arr[19,15] = 600
$arr[19,16] = 1
$arr[19,17] = 'Hortaliza'
$arr[20,0] = 11
$arr[20,1] = 'Vega Monumental Concepción'
$arr[20,2] = 'Bíobío'
$arr[20,3] = 44754
$arr[20,4] = 8
$arr[20,5] = 100112023
$arr[20,6] = 'Brócoli'
$arr[20,7] = 'Sin especificar'
$arr[20,8] = 'Primera'
$arr[20,9] = 2000
$arr[20,10] = 1000
$arr[20,11] = 1100
$arr[20,12] = 1050
$arr[20,13] = '$/unidad'
$arr[20,14] = 'Región Metropolitana'
$arr[20,15] = 1050
$arr[20,16] = 1
$arr[20,17] = 'Hortaliza'
$arr[21,0] = 11
$arr[21,1] = 'Vega Monumental Concepción'
$arr[21,2] = 'Bíobío'
$arr[21,3] = 44365
$arr[21,4] = 8
$arr[21,5] = 100112023
$arr[21,6] = 'Brócoli'
$arr[21,7] = 'Sin especificar'
$arr[21,8] = 'Primera'
$arr[21,9] = 2000
$arr[21,10] = 600
$arr[21,11] = 700
$arr[21,12] = 650
$arr[21,13] = '$/unidad'
$arr[21,14] = 'Región Metropolitana'
$arr[21,15] = 650
$arr[21,16] = 1
$arr[21,17] = 'Hortaliza'
$arr[22,0] = 11
$arr[22,1] = 'Vega Monumental Concepción'
$arr[22,2] = 'Bíobío'
$arr[22,3] = 44365
$arr[22,4] = 8
$arr[22,5] = 100112023
$arr[22,6] = 'Brócoli'
$arr[22,7] = 'Sin especificar'
$arr[22,8] = 'Segunda'
$arr[22,9] = 1000
$arr[22,10] = 500
$arr[22,11] = 500
$arr[22,12] = 500
$arr[22,13] = '$/unidad'
$arr[22,14] = 'Región Metropolitana'
$arr[22,15] = 500
$arr[22,16] = 1
$arr[22,17] = 'Hortaliza'
$arr[23,0] = 11
$arr[23,1] = 'Vega Monumental Concepción'
$arr[23,2] = 'Bíobío'
$arr[23,3] = 44736
$arr[23,4] = 8
$arr[23,5] = 100112023
$arr[23,6] = 'Brócoli'
$arr[23,7] = 'Sin especificar'
$arr[23,8] = 'Primera'
$arr[23,9] = 2000
$arr[23,10] = 700
$arr[23,11] = 750
$arr[23,12] = 725
$arr[23,13] = '$/unidad'
$arr[23,14] = 'Región Metropolitana'
$arr[23,15] = 725
$arr[23,16] = 1
$arr[23,17] = 'Hortaliza'
$arr[24,0] = 11
$arr[24,1] = 'Vega Monumental Concepción'
$arr[24,2] = 'Bíobío'
$arr[24,3] = 44736
$arr[24,4] = 8
$arr[24,5] = 100112023
$arr[24,6] = 'Brócoli'
$arr[24,7] = 'Sin especificar'
$arr[24,8] = 'Segunda'
$arr[24,9] = 1000
$arr[24,10] = 600
$arr[24,11] = 600
$arr[24,12] = 600
$arr[24,13] = '$/unidad'
$arr[24,14] = 'Región Metropolitana'
$arr[24,15] = 600
$arr[24,16] = 1
$arr[24,17] = 'Hortaliza'
$arr[25,0] = 11
$arr[25,1] = 'Vega Monumental Concepción'
$arr[25,2] = 'Bíobío'
$arr[25,3] = 44314
$arr[25,4] = 8
$arr[25,5] = 100112023
$arr[25,6] = 'Brócoli'
$arr[25,7] = 'Sin especificar'
$arr[25,8] = 'Primera'
$arr[25,9] = 1000
$arr[25,10] = 700
$arr[25,11] = 800
$arr[25,12] = 750
$arr[25,13] = '$/unidad'
$arr[25,14] = 'Región Metropolitana'
$arr[25,15] = 750
$arr[25,16] = 1
$arr[25,17] = 'Hortaliza'
$arr[26,0] = 11
$arr[26,1] = 'Vega Monumental Concepción'
$arr[26,2] = 'Bíobío'
$arr[26,3] = 44314
$arr[26,4] = 8
$arr[26,5] = 100112023
$arr[26,6] = 'Brócoli'
$arr[26,7] = 'Sin especificar'
$arr[26,8] = 'Segunda'
$arr[26,9] = 500
$arr[26,10] = 600
$arr[26,11] = 600
$arr[26,12] = 600
$arr[26,13] = '$/unidad'
$arr[26,14] = 'Región Metropolitana'
$arr[26,15] = 600
$arr[26,16] = 1
$arr[26,17] = 'Hortaliza'
$arr[27,0] = 11
$arr[27,1] = 'Vega Monumental Concepción'
$arr[27,2] = 'Bíobío'
$arr[27,3] = 44194
$arr[27,4] = 8
$arr[27,5] = 100112023
$arr[27,6] = 'Brócoli'
$arr[27,7] = 'Sin especificar'
$arr[27,8] = 'Primera'
$arr[27,9] = 1000
$arr[27,10] = 800
$arr[27,11] = 900
$arr[27,12] = 850
$arr[27,13] = '$/unidad'
$arr[27,14] = 'Región Metropolitana'
$arr[27,15] = 850
$arr[27,16] = 1
$arr[27,17] = 'Hortaliza'
$arr[28,0] = 11
$arr[28,1] = 'Vega Monumental Concepción'
$arr[28,2] = 'Bíobío'
$arr[28,3] = 44194
$arr[28,4] = 8
$arr[28,5] = 100112023
$arr[28,6] = 'Brócoli'
$arr[28,7] = 'Sin especificar'
$arr[28,8] = 'Segunda'
$arr[28,9] = 500
$arr[28,10] = 700
$arr[28,11] = 700
$arr[28,12] = 700
$arr[28,13] = '$/unidad'
$arr[28,14] = 'Región Metropolitana'
$arr[28,15] = 700
$arr[28,16] = 1
$arr[28,17] = 'Hortaliza'
$arr[29,0] = 11
$arr[29,1] = 'Vega Monumental Concepción'
$arr[29,2] = 'Bíobío'
$arr[29,3] = 44313
$arr[29,4] = 8
$arr[29,5] = 100112023
$arr[29,6] = 'Brócoli'
$arr[29,7] = 'Sin especificar'
$arr[29,8] = 'Primera'
$arr[29,9] = 2000
$arr[29,10] = 600
$arr[29,11] = 700
$arr[29,12] = 650
$arr[29,13] = '$/unidad'
$arr[29,14] = 'Región Metropolitana'
$arr[29,15] = 650
$arr[29,16] = 1
$arr[29,17] = 'Hortaliza'
$arr[30,0] = 11
$arr[30,1] = 'Vega Monumental Concepción'
$arr[30,2] = 'Bíobío'
$arr[30,3] = 44313
$arr[30,4] = 8
$arr[30,5] = 100112023
$arr[30,6] = 'Brócoli'
$arr[30,7] = 'Sin especificar'
$arr[30,8] = 'Segunda'
$arr[30,9] = 800
$arr[30,10] = 500
$arr[30,11] = 500
$arr[30,12] = 500
$arr[30,13] = '$/unidad'
$arr[30,14] = 'Región Metropolitana'
$arr[30,15] = 500
$arr[30,16] = 1
$arr[30,17] = 'Hortaliza'
$arr[31,0] = 11
$arr[31,1] = 'Vega Monumental Concepción'
$arr[31,2] = 'Bíobío'
$arr[31,3] = 44579
$arr[31,4] = 8
$arr[31,5] = 100112023
$arr[31,6] = 'Brócoli'
$arr[31,7] = 'Sin especificar'
$arr[31,8] = 'Primera'
$arr[31,9] = 2000
$arr[31,10] = 1000
$arr[31,11] = 1100
$arr[31,12] = 1050
$arr[31,13] = '$/unidad'
$arr[31,14] = 'Región Metropolitana'
$arr[31,15] = 1050
$arr[31,16] = 1
$arr[31,17] = 'Hortaliza'
$arr[32,0] = 11
$arr[32,1] = 'Vega Monumental Concepción'
$arr[32,2] = 'Bíobío'
$arr[32,3] = 44579
$arr[32,4] = 8
$arr[32,5] = 100112023
$arr[32,6] = 'Brócoli'
$arr[32,7] = 'Sin especificar'
$arr[32,8] = 'Segunda'
$arr[32,9] = 1000
$arr[32,10] = 900
$arr[32,11] = 900
$arr[32,12] = 900
$arr[32,13] = '$/unidad'
$arr[32,14] = 'Región Metropolitana'
$arr[32,15] = 900
$arr[32,16] = 1
$arr[32,17] = 'Hortaliza'
$arr[33,0] = 11
$arr[33,1] = 'Vega Monumental Concepción'
$arr[33,2] = 'Bíobío'
$arr[33,3] = 44463
$arr[33,4] = 8
$arr[33,5] = 100112023
$arr[33,6] = 'Brócoli'
$arr[33,7] = 'Sin especificar'
$arr[33,8] = 'Primera'
$arr[33,9] = 2000
$arr[33,10] = 700
$arr[33,11] = 800
$arr[33,12] = 750
$arr[33,13] = '$/unidad'
$arr[33,14] = 'Región Metropolitana'
$arr[33,15] = 750
$arr[33,16] = 1
$arr[33,17] = 'Hortaliza'
$arr[34,0] = 11
$arr[34,1] = 'Vega Monumental Concepción'
$arr[34,2] = 'Bíobío'
$arr[34,3] = 44463
$arr[34,4] = 8
$arr[34,5] = 100112023
$arr[34,6] = 'Brócoli'
$arr[34,7] = 'Sin especificar'
$arr[34,8] = 'Segunda'
$arr[34,9] = 1000
$arr[34,10] = 600
$arr[34,11] = 600
$arr[34,12] = 600
$arr[34,13] = '$/unidad'
$arr[34,14] = 'Región Metropolitana'
$arr[34,15] = 600
$arr[34,16] = 1
$arr[34,17] = 'Hortaliza'
$arr[35,0] = 11
$arr[35,1] = 'Vega Monumental Concepción'
$arr[35,2] = 'Bíobío'
$arr[35,3] = 44685
$arr[35,4] = 8
$arr[35,5] = 100112023
$arr[35,6] = 'Brócoli'
$arr[35,7] = 'Sin especificar'
$arr[35,8] = 'Primera'
$arr[35,9] = 2000
$arr[35,10] = 700
$arr[35,11] = 800
$arr[35,12] = 750
$arr[35,13] = '$/unidad'
$arr[35,14] = 'Región Metropolitana'
$arr[35,15] = 750
$arr[35,16] = 1
$arr[35,17] = 'Hortaliza'
$arr[36,0] = 11
$arr[36,1] = 'Vega Monumental Concepción'
$arr[36,2] = 'Bíobío'
$arr[36,3] = 44685
$arr[36,4] = 8
$arr[36,5] = 100112023
$arr[36,6] = 'Brócoli'
$arr[36,7] = 'Sin especificar'
$arr[36,8] = 'Segunda'
$arr[36,9] = 1000
$arr[36,10] = 500
$arr[36,11] = 500
$arr[36,12] = 500
$arr[36,13] = '$/unidad'
$arr[36,14] = 'Región Metropolitana'
$arr[36,15] = 500
$arr[36,16] = 1
$arr[36,17] = 'Hortaliza'
$arr[37,0] = 11
$arr[37,1] = 'Vega Monumental Concepción'
$arr[37,2] = 'Bíobío'
$arr[37,3] = 44327
$arr[37,4] = 8
$arr[37,5] = 100112023
$arr[37,6] = 'Brócoli'
$arr[37,7] = 'Sin especificar'
$arr[37,8] = 'Primera'
$arr[37,9] = 2000
$arr[37,10] = 600
$arr[37,11] = 700
$arr[37,12] = 650
$arr[37,13] = '$/unidad'
$arr[37,14] = 'Región Metropolitana'
$arr[37,15] = 650
$arr[37,16] = 1
$arr[37,17] = 'Hortaliza'
$arr[38,0] = 11
$arr[38,1] = 'Vega Monumental Concepción'
$arr[38,2] = 'Bíobío'
$arr[38,3] = 44327
$arr[38,4] = 8
$arr[38,5] = 100112023
$arr[38,6] = 'Brócoli'
$arr[38,7] = 'Sin especificar'
$arr[38,8] = 'Segunda'
$arr[38,9] = 1000
$arr[38,10] = 500
$arr[38,11] = 500
$arr[38,12] = 500
$arr[38,13] = '$/unidad'
$arr[38,14] = 'Región Metropolitana'
$arr[38,15] = 500
$arr[38,16] = 1
$arr[38,17] = 'Hortaliza'
$arr[39,0] = 11
$arr[39,1] = 'Vega Monumental Concepción'
$arr[39,2] = 'Bíobío'
$arr[39,3] = 44673
$arr[39,4] = 8
$arr[39,5] = 100112023
$arr[39,6] = 'Brócoli'
$arr[39,7] = 'Sin especificar'
$arr[39,8] = 'Primera'
$arr[39,9] = 2500
$arr[39,10] = 750
$arr[39,11] = 800
$arr[39,12] = 780
$arr[39,13] = '$/unidad'
$arr[39,14] = 'Región Metropolitana'
$arr[39,15] = 780
$arr[39,16] = 1
$arr[39,17] = 'Hortaliza'
$arr[40,0] = 11
$arr[40,1] = 'Vega Monumental Concepción'
$arr[40,2] = 'Bíobío'
$arr[40,3] = 44460
$arr[40,4] = 8
$arr[40,5] = 100112023
$arr[40,6] = 'Brócoli'
$arr[40,7] = 'Sin especificar'
$arr[40,8] = 'Primera'
$arr[40,9] = 1000
$arr[40,10] = 600
$arr[40,11] = 700
$arr[40,12] = 650
$arr[40,13] = '$/unidad'
$arr[40,14] = 'Región Metropolitana'
$arr[40,15] = 650
$arr[40,16] = 1
$arr[40,17] = 'Hortaliza'
$arr[41,0] = 11
$arr[41,1] = 'Vega Monumental Concepción'
$arr[41,2] = 'Bíobío'
$arr[41,3] = 44460
$arr[41,4] = 8
$arr[41,5] = 100112023
$arr[41,6] = 'Brócoli'
$arr[41,7] = 'Sin especificar'
$arr[41,8] = 'Segunda'
$arr[41,9] = 500
$arr[41,10] = 500
$arr[41,11] = 500
$arr[41,12] = 500
$arr[41,13] = '$/unidad'
$arr[41,14] = 'Región Metropolitana'
$arr[41,15] = 500
$arr[41,16] = 1
$arr[41,17] = 'Hortaliza'
$arr[42,0] = 11
$arr[42,1] = 'Vega Monumental Concepción'
$arr[42,2] = 'Bíobío'
$arr[42,3] = 44649
$arr[42,4] = 8
$arr[42,5] = 100112023
$arr[42,6] = 'Brócoli'
$arr[42,7] = 'Sin especificar'
$arr[42,8] = 'Primera'
$arr[42,9] = 1500
$arr[42,10] = 1000
$arr[42,11] = 1000
$arr[42,12] = 1000
$arr[42,13] = '$/unidad'
$arr[42,14] = 'Región Metropolitana'
$arr[42,15] = 1000
$arr[42,16] = 1
$arr[42,17] = 'Hortaliza'
$arr[43,0] = 11
$arr[43,1] = 'Vega Monumental Concepción'
$arr[43,2] = 'Bíobío'
$arr[43,3] = 44649
$arr[43,4] = 8
$arr[43,5] = 100112023
$arr[43,6] = 'Brócoli'
$arr[43,7] = 'Sin especificar'
$arr[43,8] = 'Segunda'
$arr[43,9] = 1500
$arr[43,10] = 800
$arr[43,11] = 800
$arr[43,12] = 800
$arr[43,13] = '$/unidad'
$arr[43,14] = 'Región Metropolitana'
$arr[43,15] = 800
$arr[43,16] = 1
$arr[43,17] = 'Hortaliza'
$arr[44,0] = 11
$arr[44,1] = 'Vega Monumental Concepción'
$arr[44,2] = 'Bíobío'
$arr[44,3] = 44607
$arr[44,4] = 8
$arr[44,5] = 100112023
$arr[44,6] = 'Brócoli'
$arr[44,7] = 'Sin especificar'
$arr[44,8] = 'Primera'
$arr[44,9] = 2000
$arr[44,10] = 900
$arr[44,11] = 1000
$arr[44,12] = 950
$arr[44,13] = '$/unidad'
$arr[44,14] = 'Provincia de Chacabuco'
$arr[44,15] = 950
$arr[44,16] = 1
$arr[44,17] = 'Hortaliza'
$arr[45,0] = 11
$arr[45,1] = 'Vega Monumental Concepción'
$arr[45,2] = 'Bíobío'
$arr[45,3] = 44607
$arr[45,4] = 8
$arr[45,5] = 100112023
$arr[45,6] = 'Brócoli'
$arr[45,7] = 'Sin especificar'
$arr[45,8] = 'Segunda'
$arr[45,9] = 1000
$arr[45,10] = 800
$arr[45,11] = 800
$arr[45,12] = 800
$arr[45,13] = '$/unidad'
$arr[45,14] = 'Provincia de Chacabuco'
$arr[45,15] = 800
$arr[45,16] = 1
$arr[45,17] = 'Hortaliza'
$arr[46,0] = 11
$arr[46,1] = 'Vega Monumental Concepción'
$arr[46,2] = 'Bíobío'
$arr[46,3] = 44292
$arr[46,4] = 8
$arr[46,5] = 100112023
$arr[46,6] = 'Brócoli'
$arr[46,7] = 'Sin especificar'
$arr[46,8] = 'Primera'
$arr[46,9] = 1000
$arr[46,10] = 700
$arr[46,11] = 800
$arr[46,12] = 750
$arr[46,13] = '$/unidad'
$arr[46,14] = 'Región Metropolitana'
$arr[46,15] = 750
$arr[46,16] = 1
$arr[46,17] = 'Hortaliza'
$arr[47,0] = 11
$arr[47,1] = 'Vega Monumental Concepción'
$arr[47,2] = 'Bíobío'
$arr[47,3] = 44292
$arr[47,4] = 8
$arr[47,5] = 100112023
$arr[47,6] = 'Brócoli'
$arr[47,7] = 'Sin especificar'
$arr[47,8] = 'Segunda'
$arr[47,9] = 500
$arr[47,10] = 600
$arr[47,11] = 600
$arr[47,12] = 600
$arr[47,13] = '$/unidad'
$arr[47,14] = 'Región Metropolitana'
$arr[47,15] = 600
$arr[47,16] = 1
$arr[47,17] = 'Hortaliza'
$arr[48,0] = 11
$arr[48,1] = 'Vega Monumental Concepción'
$arr[48,2] = 'Bíobío'
$arr[48,3] = 44348
$arr[48,4] = 8
$arr[48,5] = 100112023
$arr[48,6] = 'Brócoli'
$arr[48,7] = 'Sin especificar'
$arr[48,8] = 'Primera'
$arr[48,9] = 2000
$arr[48,10] = 600
$arr[48,11] = 700
$arr[48,12] = 650
$arr[48,13] = '$/unidad'
$arr[48,14] = 'Región Metropolitana'
$arr[48,15] = 650
$arr[48,16] = 1
$arr[48,17] = 'Hortaliza'
$arr[49,0] = 11
$arr[49,1] = 'Vega Monumental Concepción'
$arr[49,2] = 'Bíobío'
$arr[49,3] = 44348
$arr[49,4] = 8
$arr[49,5] = 100112023
$arr[49,6] = 'Brócoli'
$arr[49,7] = 'Sin especificar'
$arr[49,8] = 'Segunda'
$arr[49,9] = 1000
$arr[49,10] = 500
$arr[49,11] = 500
$arr[49,12] = 500
$arr[49,13] = '$/unidad'
$arr[49,14] = 'Región Metropolitana'
$arr[49,15] = 500
$arr[49,16] = 1
$arr[49,17] = 'Hortaliza'
$arr[50,0] = 11
$arr[50,1] = 'Vega Monumental Concepción'
$arr[50,2] = 'Bíobío'
$arr[50,3] = 44341
$arr[50,4] = 8
$arr[50,5] = 100112023
$arr[50,6] = 'Brócoli'
$arr[50,7] = 'Sin especificar'
$arr[50,8] = 'Primera'
$arr[50,9] = 2000
$arr[50,10] = 600
$arr[50,11] = 650
$arr[50,12] = 625
$arr[50,13] = '$/unidad'
$arr[50,14] = 'Región Metropolitana'
$arr[50,15] = 625
$arr[50,16] = 1
$arr[50,17] = 'Hortaliza'
$arr[51,0] = 11
$arr[51,1] = 'Vega Monumental Concepción'
$arr[51,2] = 'Bíobío'
$arr[51,3] = 44341
$arr[51,4] = 8
$arr[51,5] = 100112023
$arr[51,6] = 'Brócoli'
$arr[51,7] = 'Sin especificar'
$arr[51,8] = 'Segunda'
$arr[51,9] = 1000
$arr[51,10] = 500
$arr[51,11] = 500
$arr[51,12] = 500
$arr[51,13] = '$/unidad'
$arr[51,14] = 'Región Metropolitana'
$arr[51,15] = 500
$arr[51,16] = 1
$arr[51,17] = 'Hortaliza'
$arr[52,0] = 11
$arr[52,1] = 'Vega Monumental Concepción'
$arr[52,2] = 'Bíobío'
$arr[52,3] = 44370
$arr[52,4] = 8
$arr[52,5] = 100112023
$arr[52,6] = 'Brócoli'
$arr[52,7] = 'Sin especificar'
$arr[52,8] = 'Primera'
$arr[52,9] = 2000
$arr[52,10] = 600
$arr[52,11] = 700
$arr[52,12] = 650
$arr[52,13] = '$/unidad'
$arr[52,14] = 'Región del Maule'
$arr[52,15] = 650
$arr[52,16] = 1
$arr[52,17] = 'Hortaliza'
$arr[53,0] = 11
$arr[53,1] = 'Vega Monumental Concepción'
$arr[53,2] = 'Bíobío'
$arr[53,3] = 44370
$arr[53,4] = 8
$arr[53,5] = 100112023
$arr[53,6] = 'Brócoli'
$arr[53,7] = 'Sin especificar'
$arr[53,8] = 'Segunda'
$arr[53,9] = 1000
$arr[53,10] = 500
$arr[53,11] = 500
$arr[53,12] = 500
$arr[53,13] = '$/unidad'
$arr[53,14] = 'Región del Maule'
$arr[53,15] = 500
$arr[53,16] = 1
$arr[53,17] = 'Hortaliza'
$arr[54,0] = 11
$arr[54,1] = 'Vega Monumental Concepción'
$arr[54,2] = 'Bíobío'
$arr[54,3] = 44622
$arr[54,4] = 8
$arr[54,5] = 100112023
$arr[54,6] = 'Brócoli'
$arr[54,7] = 'Sin especificar'
$arr[54,8] = 'Primera'
$arr[54,9] = 1500
$arr[54,10] = 900
$arr[54,11] = 1000
$arr[54,12] = 953
$arr[54,13] = '$/unidad'
$arr[54,14] = 'Región Metropolitana'
$arr[54,15] = 953
$arr[54,16] = 1
$arr[54,17] = 'Hortaliza'
$arr[55,0] = 11
$arr[55,1] = 'Vega Monumental Concepción'
$arr[55,2] = 'Bíobío'
$arr[55,3] = 44565
$arr[55,4] = 8
$arr[55,5] = 100112023
$arr[55,6] = 'Brócoli'
$arr[55,7] = 'Sin especificar'
$arr[55,8] = 'Primera'
$arr[55,9] = 2000
$arr[55,10] = 700
$arr[55,11] = 800
$arr[55,12] = 750
$arr[55,13] = '$/unidad'
$arr[55,14] = 'Región Metropolitana'
$arr[55,15] = 750
$arr[55,16] = 1
$arr[55,17] = 'Hortaliza'
$arr[56,0] = 11
$arr[56,1] = 'Vega Monumental Concepción'
$arr[56,2] = 'Bíobío'
$arr[56,3] = 44565
$arr[56,4] = 8
$arr[56,5] = 100112023
$arr[56,6] = 'Brócoli'
$arr[56,7] = 'Sin especificar'
$arr[56,8] = 'Segunda'
$arr[56,9] = 1000
$arr[56,10] = 600
$arr[56,11] = 600
$arr[56,12] = 600
$arr[56,13] = '$/unidad'
$arr[56,14] = 'Región Metropolitana'
$arr[56,15] = 600
$arr[56,16] = 1
$arr[56,17] = 'Hortaliza'
$arr[57,0] = 11
$arr[57,1] = 'Vega Monumental Concepción'
$arr[57,2] = 'Bíobío'
$arr[57,3] = 44748
$arr[57,4] = 8
$arr[57,5] = 100112023
$arr[57,6] = 'Brócoli'
$arr[57,7] = 'Sin especificar'
$arr[57,8] = 'Primera'
$arr[57,9] = 650
$arr[57,10] = 700
$arr[57,11] = 750
$arr[57,12] = 723
$arr[57,13] = '$/unidad'
$arr[57,14] = 'Región Metropolitana'
$arr[57,15] = 723
$arr[57,16] = 1
$arr[57,17] = 'Hortaliza'
$arr[58,0] = 11
$arr[58,1] = 'Vega Monumental Concepción'
$arr[58,2] = 'Bíobío'
$arr[58,3] = 44259
$arr[58,4] = 8
$arr[58,5] = 100112023
$arr[58,6] = 'Brócoli'
$arr[58,7] = 'Sin especificar'
$arr[58,8] = 'Primera'
$arr[58,9] = 500
$arr[58,10] = 1000
$arr[58,11] = 1000
$arr[58,12] = 1000
$arr[58,13] = '$/unidad'
$arr[58,14] = 'Región Metropolitana'
$arr[58,15] = 1000
$arr[58,16] = 1
$arr[58,17] = 'Hortaliza'
$arr[59,0] = 11
$arr[59,1] = 'Vega Monumental Concepción'
$arr[59,2] = 'Bíobío'
$arr[59,3] = 44259
$arr[59,4] = 8
$arr[59,5] = 100112023
$arr[59,6] = 'Brócoli'
$arr[59,7] = 'Sin especificar'
$arr[59,8] = 'Segunda'
$arr[59,9] = 500
$arr[59,10] = 700
$arr[59,11] = 700
$arr[59,12] = 700
$arr[59,13] = '$/unidad'
$arr[59,14] = 'Región Metropolitana'
$arr[59,15] = 700
$arr[59,16] = 1
$arr[59,17] = 'Hortaliza'
$arr[60,0] = 11
$arr[60,1] = 'Vega Monumental Concepción'
$arr[60,2] = 'Bíobío'
$arr[60,3] = 44202
$arr[60,4] = 8
$arr[60,5] = 100112023
$arr[60,6] = 'Brócoli'
$arr[60,7] = 'Sin especificar'
$arr[60,8] = 'Primera'
$arr[60,9] = 800
$arr[60,10] = 600
$arr[60,11] = 700
$arr[60,12] = 650
$arr[60,13] = '$/unidad'
$arr[60,14] = 'Región del Maule'
$arr[60,15] = 650
$arr[60,16] = 1
$arr[60,17] = 'Hortaliza'
$arr[61,0] = 11
$arr[61,1] = 'Vega Monumental Concepción'
$arr[61,2] = 'Bíobío'
$arr[61,3] = 44202
$arr[61,4] = 8
$arr[61,5] = 100112023
$arr[61,6] = 'Brócoli'
$arr[61,7] = 'Sin especificar'
$arr[61,8] = 'Segunda'
$arr[61,9] = 400
$arr[61,10] = 500
$arr[61,11] = 500
$arr[61,12] = 500
$arr[61,13] = '$/unidad'
$arr[61,14] = 'Región del Maule'
$arr[61,15] = 500
$arr[61,16] = 1
$arr[61,17] = 'Hortaliza'
$arr[62,0] = 11
$arr[62,1] = 'Vega Monumental Concepción'
$arr[62,2] = 'Bíobío'
$arr[62,3] = 44162
$arr[62,4] = 8
$arr[62,5] = 100112023
$arr[62,6] = 'Brócoli'
$arr[62,7] = 'Sin especificar'
$arr[62,8] = 'Primera'
$arr[62,9] = 800
$arr[62,10] = 600
$arr[62,11] = 700
$arr[62,12] = 650
$arr[62,13] = '$/unidad'
$arr[62,14] = 'Región del Maule'
$arr[62,15] = 650
$arr[62,16] = 1
$arr[62,17] = 'Hortaliza'
$arr[63,0] = 11
$arr[63,1] = 'Vega Monumental Concepción'
$arr[63,2] = 'Bíobío'
$arr[63,3] = 44162
$arr[63,4] = 8
$arr[63,5] = 100112023
$arr[63,6] = 'Brócoli'
$arr[63,7] = 'Sin especificar'
$arr[63,8] = 'Segunda'
$arr[63,9] = 400
$arr[63,10] = 500
$arr[63,11] = 500
$arr[63,12] = 500
$arr[63,13] = '$/unidad'
$arr[63,14] = 'Región del Maule'
$arr[63,15] = 500
$arr[63,16] = 1
$arr[63,17] = 'Hortaliza'
$arr[64,0] = 11
$arr[64,1] = 'Vega Monumental Concepción'
$arr[64,2] = 'Bíobío'
$arr[64,3] = 44726
$arr[64,4] = 8
$arr[64,5] = 100112023
$arr[64,6] = 'Brócoli'
$arr[64,7] = 'Sin especificar'
$arr[64,8] = 'Primera'
$arr[64,9] = 2200
$arr[64,10] = 800
$arr[64,11] = 900
$arr[64,12] = 855
$arr[64,13] = '$/unidad'
$arr[64,14] = 'Región Metropolitana'
$arr[64,15] = 855
$arr[64,16] = 1
$arr[64,17] = 'Hortaliza'
$arr[65,0] = 11
$arr[65,1] = 'Vega Monumental Concepción'
$arr[65,2] = 'Bíobío'
$arr[65,3] = 44427
$arr[65,4] = 8
$arr[65,5] = 100112023
$arr[65,6] = 'Brócoli'
$arr[65,7] = 'Sin especificar'
$arr[65,8] = 'Primera'
$arr[65,9] = 2000
$arr[65,10] = 600
$arr[65,11] = 700
$arr[65,12] = 650
$arr[65,13] = '$/unidad'
$arr[65,14] = 'Región Metropolitana'
$arr[65,15] = 650
$arr[65,16] = 1
$arr[65,17] = 'Hortaliza'
$arr[66,0] = 11
$arr[66,1] = 'Vega Monumental Concepción'
$arr[66,2] = 'Bíobío'
$arr[66,3] = 44427
$arr[66,4] = 8
$arr[66,5] = 100112023
$arr[66,6] = 'Brócoli'
$arr[66,7] = 'Sin especificar'
$arr[66,8] = 'Segunda'
$arr[66,9] = 1000
$arr[66,10] = 500
$arr[66,11] = 500
$arr[66,12] = 500
$arr[66,13] = '$/unidad'
$arr[66,14] = 'Región Metropolitana'
$arr[66,15] = 500
$arr[66,16] = 1
$arr[66,17] = 'Hortaliza'
$arr[67,0] = 11
$arr[67,1] = 'Vega Monumental Concepción'
$arr[67,2] = 'Bíobío'
$arr[67,3] = 44441
$arr[67,4] = 8
$arr[67,5] = 100112023
$arr[67,6] = 'Brócoli'
$arr[67,7] = 'Sin especificar'
$arr[67,8] = 'Primera'
$arr[67,9] = 1000
$arr[67,10] = 700
$arr[67,11] = 800
$arr[67,12] = 750
$arr[67,13] = '$/unidad'
$arr[67,14] = 'Región Metropolitana'
$arr[67,15] = 750
$arr[67,16] = 1
$arr[67,17] = 'Hortaliza'
$arr[68,0] = 11
$arr[68,1] = 'Vega Monumental Concepción'
$arr[68,2] = 'Bíobío'
$arr[68,3] = 44441
$arr[68,4] = 8
$arr[68,5] = 100112023
$arr[68,6] = 'Brócoli'
$arr[68,7] = 'Sin especificar'
$arr[68,8] = 'Segunda'
$arr[68,9] = 500
$arr[68,10] = 600
$arr[68,11] = 600
$arr[68,12] = 600
$arr[68,13] = '$/unidad'
$arr[68,14] = 'Región Metropolitana'
$arr[68,15] = 600
$arr[68,16] = 1
$arr[68,17] = 'Hortaliza'
$arr[69,0] = 11
$arr[69,1] = 'Vega Monumental Concepción'
$arr[69,2] = 'Bíobío'
$arr[69,3] = 44708
$arr[69,4] = 8
$arr[69,5] = 100112023
$arr[69,6] = 'Brócoli'
$arr[69,7] = 'Sin especificar'
$arr[69,8] = 'Primera'
$arr[69,9] = 4000
$arr[69,10] = 700
$arr[69,11] = 800
$arr[69,12] = 750
$arr[69,13] = '$/unidad'
$arr[69,14] = 'Provincia de Chacabuco'
$arr[69,15] = 750
$arr[69,16] = 1
$arr[69,17] = 'Hortaliza'
$arr[70,0] = 11
$arr[70,1] = 'Vega Monumental Concepción'
$arr[70,2] = 'Bíobío'
$arr[70,3] = 44708
$arr[70,4] = 8
$arr[70,5] = 100112023
$arr[70,6] = 'Brócoli'
$arr[70,7] = 'Sin especificar'
$arr[70,8] = 'Segunda'
$arr[70,9] = 2000
$arr[70,10] = 600
$arr[70,11] = 600
$arr[70,12] = 600
$arr[70,13] = '$/unidad'
$arr[70,14] = 'Provincia de Chacabuco'
$arr[70,15] = 600
$arr[70,16] = 1
$arr[70,17] = 'Hortaliza'
$arr[71,0] = 11
$arr[71,1] = 'Vega Monumental Concepción'
$arr[71,2] = 'Bíobío'
$arr[71,3] = 44518
$arr[71,4] = 8
$arr[71,5] = 100112023
$arr[71,6] = 'Brócoli'
$arr[71,7] = 'Sin especificar'
$arr[71,8] = 'Primera'
$arr[71,9] = 1100
$arr[71,10] = 700
$arr[71,11] = 800
$arr[71,12] = 745
$arr[71,13] = '$/unidad'
$arr[71,14] = 'Región Metropolitana'
$arr[71,15] = 745
$arr[71,16] = 1
$arr[71,17] = 'Hortaliza'
$arr[72,0] = 11
$arr[72,1] = 'Vega Monumental Concepción'
$arr[72,2] = 'Bíobío'
$arr[72,3] = 44628
$arr[72,4] = 8
$arr[72,5] = 100112023
$arr[72,6] = 'Brócoli'
$arr[72,7] = 'Sin especificar'
$arr[72,8] = 'Primera'
$arr[72,9] = 1500
$arr[72,10] = 1000
$arr[72,11] = 1000
$arr[72,12] = 1000
$arr[72,13] = '$/unidad'
$arr[72,14] = 'Región Metropolitana'
$arr[72,15] = 1000
$arr[72,16] = 1
$arr[72,17] = 'Hortaliza'
$arr[73,0] = 11
$arr[73,1] = 'Vega Monumental Concepción'
$arr[73,2] = 'Bíobío'
$arr[73,3] = 44628
$arr[73,4] = 8
$arr[73,5] = 100112023
$arr[73,6] = 'Brócoli'
$arr[73,7] = 'Sin especificar'
$arr[73,8] = 'Segunda'
$arr[73,9] = 1200
$arr[73,10] = 800
$arr[73,11] = 800
$arr[73,12] = 800
$arr[73,13] = '$/unidad'
$arr[73,14] = 'Región Metropolitana'
$arr[73,15] = 800
$arr[73,16] = 1
$arr[73,17] = 'Hortaliza'
$arr[74,0] = 11
$arr[74,1] = 'Vega Monumental Concepción'
$arr[74,2] = 'Bíobío'
$arr[74,3] = 44483
$arr[74,4] = 8
$arr[74,5] = 100112023
$arr[74,6] = 'Brócoli'
$arr[74,7] = 'Sin especificar'
$arr[74,8] = 'Primera'
$arr[74,9] = 2100
$arr[74,10] = 600
$arr[74,11] = 650
$arr[74,12] = 629
$arr[74,13] = '$/unidad'
$arr[74,14] = 'Provincia de Melipilla'
$arr[74,15] = 629
$arr[74,16] = 1
$arr[74,17] = 'Hortaliza'
$arr[75,0] = 11
$arr[75,1] = 'Vega Monumental Concepción'
$arr[75,2] = 'Bíobío'
$arr[75,3] = 44469
$arr[75,4] = 8
$arr[75,5] = 100112023
$arr[75,6] = 'Brócoli'
$arr[75,7] = 'Sin especificar'
$arr[75,8] = 'Primera'
$arr[75,9] = 1000
$arr[75,10] = 700
$arr[75,11] = 800
$arr[75,12] = 750
$arr[75,13] = '$/unidad'
$arr[75,14] = 'Región Metropolitana'
$arr[75,15] = 750
$arr[75,16] = 1
$arr[75,17] = 'Hortaliza'
$arr[76,0] = 11
$arr[76,1] = 'Vega Monumental Concepción'
$arr[76,2] = 'Bíobío'
$arr[76,3] = 44469
$arr[76,4] = 8
$arr[76,5] = 100112023
$arr[76,6] = 'Brócoli'
$arr[76,7] = 'Sin especificar'
$arr[76,8] = 'Segunda'
$arr[76,9] = 500
$arr[76,10] = 600
$arr[76,11] = 600
$arr[76,12] = 600
$arr[76,13] = '$/unidad'
$arr[76,14] = 'Región Metropolitana'
$arr[76,15] = 600
$arr[76,16] = 1
$arr[76,17] = 'Hortaliza'
$arr[77,0] = 11
$arr[77,1] = 'Vega Monumental Concepción'
$arr[77,2] = 'Bíobío'
$arr[77,3] = 44776
$arr[77,4] = 8
$arr[77,5] = 100112023
$arr[77,6] = 'Brócoli'
$arr[77,7] = 'Sin especificar'
$arr[77,8] = 'Primera'
$arr[77,9] = 1000
$arr[77,10] = 900
$arr[77,11] = 1000
$arr[77,12] = 950
$arr[77,13] = '$/unidad'
$arr[77,14] = 'Región Metropolitana'
$arr[77,15] = 950
$arr[77,16] = 1
$arr[77,17] = 'Hortaliza'
$arr[78,0] = 11
$arr[78,1] = 'Vega Monumental Concepción'
$arr[78,2] = 'Bíobío'
$arr[78,3] = 44776
$arr[78,4] = 8
$arr[78,5] = 100112023
$arr[78,6] = 'Brócoli'
$arr[78,7] = 'Sin especificar'
$arr[78,8] = 'Segunda'
$arr[78,9] = 1000
$arr[78,10] = 700
$arr[78,11] = 700
$arr[78,12] = 700
$arr[78,13] = '$/unidad'
$arr[78,14] = 'Región Metropolitana'
$arr[78,15] = 700
$arr[78,16] = 1
$arr[78,17] = 'Hortaliza'
$arr[79,0] = 11
$arr[79,1] = 'Vega Monumental Concepción'
$arr[79,2] = 'Bíobío'
$arr[79,3] = 44204
$arr[79,4] = 8
$arr[79,5] = 100112023
$arr[79,6] = 'Brócoli'
$arr[79,7] = 'Sin especificar'
$arr[79,8] = 'Primera'
$arr[79,9] = 800
$arr[79,10] = 700
$arr[79,11] = 800
$arr[79,12] = 750
$arr[79,13] = '$/unidad'
$arr[79,14] = 'Región del Maule'
$arr[79,15] = 750
$arr[79,16] = 1
$arr[79,17] = 'Hortaliza'
$arr[80,0] = 11
$arr[80,1] = 'Vega Monumental Concepción'
$arr[80,2] = 'Bíobío'
$arr[80,3] = 44204
$arr[80,4] = 8
$arr[80,5] = 100112023
$arr[80,6] = 'Brócoli'
$arr[80,7] = 'Sin especificar'
$arr[80,8] = 'Segunda'
$arr[80,9] = 400
$arr[80,10] = 600
$arr[80,11] = 600
$arr[80,12] = 600
$arr[80,13] = '$/unidad'
$arr[80,14] = 'Región del Maule'
$arr[80,15] = 600
$arr[80,16] = 1
$arr[80,17] = 'Hortaliza'
$arr[81,0] = 11
$arr[81,1] = 'Vega Monumental Concepción'
$arr[81,2] = 'Bíobío'
$arr[81,3] = 44484
$arr[81,4] = 8
$arr[81,5] = 100112023
$arr[81,6] = 'Brócoli'
$arr[81,7] = 'Sin especificar'
$arr[81,8] = 'Primera'
$arr[81,9] = 2700
$arr[81,10] = 600
$arr[81,11] = 650
$arr[81,12] = 622
$arr[81,13] = '$/unidad'
$arr[81,14] = 'Región Metropolitana'
$arr[81,15] = 622
$arr[81,16] = 1
$arr[81,17] = 'Hortaliza'
$arr[82,0] = 11
$arr[82,1] = 'Vega Monumental Concepción'
$arr[82,2] = 'Bíobío'
$arr[82,3] = 44484
$arr[82,4] = 8
$arr[82,5] = 100112023
$arr[82,6] = 'Brócoli'
$arr[82,7] = 'Sin especificar'
$arr[82,8] = 'Segunda'
$arr[82,9] = 1350
$arr[82,10] = 400
$arr[82,11] = 500
$arr[82,12] = 448
$arr[82,13] = '$/unidad'
$arr[82,14] = 'Región Metropolitana'
$arr[82,15] = 448
$arr[82,16] = 1
$arr[82,17] = 'Hortaliza'
$arr[83,0] = 11
$arr[83,1] = 'Vega Monumental Concepción'
$arr[83,2] = 'Bíobío'
$arr[83,3] = 44266
$arr[83,4] = 8
$arr[83,5] = 100112023
$arr[83,6] = 'Brócoli'
$arr[83,7] = 'Sin especificar'
$arr[83,8] = 'Primera'
$arr[83,9] = 800
$arr[83,10] = 700
$arr[83,11] = 750
$arr[83,12] = 725
$arr[83,13] = '$/unidad'
$arr[83,14] = 'Región del Maule'
$arr[83,15] = 725
$arr[83,16] = 1
$arr[83,17] = 'Hortaliza'
$arr[84,0] = 11
$arr[84,1] = 'Vega Monumental Concepción'
$arr[84,2] = 'Bíobío'
$arr[84,3] = 44266
$arr[84,4] = 8
$arr[84,5] = 100112023
$arr[84,6] = 'Brócoli'
$arr[84,7] = 'Sin especificar'
$arr[84,8] = 'Segunda'
$arr[84,9] = 400
$arr[84,10] = 600
$arr[84,11] = 600
$arr[84,12] = 600
$arr[84,13] = '$/unidad'
$arr[84,14] = 'Región del Maule'
$arr[84,15] = 600
$arr[84,16] = 1
$arr[84,17] = 'Hortaliza'
$arr[85,0] = 11
$arr[85,1] = 'Vega Monumental Concepción'
$arr[85,2] = 'Bíobío'
$arr[85,3] = 44229
$arr[85,4] = 8
$arr[85,5] = 100112023
$arr[85,6] = 'Brócoli'
$arr[85,7] = 'Sin especificar'
$arr[85,8] = 'Primera'
$arr[85,9] = 1000
$arr[85,10] = 600
$arr[85,11] = 700
$arr[85,12] = 650
$arr[85,13] = '$/unidad'
$arr[85,14] = 'Región Metropolitana'
$arr[85,15] = 650
$arr[85,16] = 1
$arr[85,17] = 'Hortaliza'
$arr[86,0] = 11
$arr[86,1] = 'Vega Monumental Concepción'
$arr[86,2] = 'Bíobío'
$arr[86,3] = 44229
$arr[86,4] = 8
$arr[86,5] = 100112023
$arr[86,6] = 'Brócoli'
$arr[86,7] = 'Sin especificar'
$arr[86,8] = 'Segunda'
$arr[86,9] = 500
$arr[86,10] = 500
$arr[86,11] = 500
$arr[86,12] = 500
$arr[86,13] = '$/unidad'
$arr[86,14] = 'Región Metropolitana'
$arr[86,15] = 500
$arr[86,16] = 1
$arr[86,17] = 'Hortaliza'
$arr[87,0] = 11
$arr[87,1] = 'Vega Monumental Concepción'
$arr[87,2] = 'Bíobío'
$arr[87,3] = 44336
$arr[87,4] = 8
$arr[87,5] = 100112023
$arr[87,6] = 'Brócoli'
$arr[87,7] = 'Sin especificar'
$arr[87,8] = 'Primera'
$arr[87,9] = 2000
$arr[87,10] = 600
$arr[87,11] = 700
$arr[87,12] = 650
$arr[87,13] = '$/unidad'
$arr[87,14] = 'Región Metropolitana'
$arr[87,15] = 650
$arr[87,16] = 1
$arr[87,17] = 'Hortaliza'
$arr[88,0] = 11
$arr[88,1] = 'Vega Monumental Concepción'
$arr[88,2] = 'Bíobío'
$arr[88,3] = 44336
$arr[88,4] = 8
$arr[88,5] = 100112023
$arr[88,6] = 'Brócoli'
$arr[88,7] = 'Sin especificar'
$arr[88,8] = 'Segunda'
$arr[88,9] = 1000
$arr[88,10] = 500
$arr[88,11] = 500
$arr[88,12] = 500
$arr[88,13] = '$/unidad'
$arr[88,14] = 'Región Metropolitana'
$arr[88,15] = 500
$arr[88,16] = 1
$arr[88,17] = 'Hortaliza'
$arr[89,0] = 11
$arr[89,1] = 'Vega Monumental Concepción'
$arr[89,2] = 'Bíobío'
$arr[89,3] = 44488
$arr[89,4] = 8
$arr[89,5] = 100112023
$arr[89,6] = 'Brócoli'
$arr[89,7] = 'Sin especificar'
$arr[89,8] = 'Primera'
$arr[89,9] = 1000
$arr[89,10] = 700
$arr[89,11] = 800
$arr[89,12] = 750
$arr[89,13] = '$/unidad'
$arr[89,14] = 'Región Metropolitana'
$arr[89,15] = 750
$arr[89,16] = 1
$arr[89,17] = 'Hortaliza'
$arr[90,0] = 11
$arr[90,1] = 'Vega Monumental Concepción'
$arr[90,2] = 'Bíobío'
$arr[90,3] = 44488
$arr[90,4] = 8
$arr[90,5] = 100112023
$arr[90,6] = 'Brócoli'
$arr[90,7] = 'Sin especificar'
$arr[90,8] = 'Segunda'
$arr[90,9] = 500
$arr[90,10] = 600
$arr[90,11] = 600
$arr[90,12] = 600
$arr[90,13] = '$/unidad'
$arr[90,14] = 'Región Metropolitana'
$arr[90,15] = 600
$arr[90,16] = 1
$arr[90,17] = 'Hortaliza'
$arr[91,0] = 11
$arr[91,1] = 'Vega Monumental Concepción'
$arr[91,2] = 'Bíobío'
$arr[91,3] = 44196
$arr[91,4] = 8
$arr[91,5] = 100112023
$arr[91,6] = 'Brócoli'
$arr[91,7] = 'Sin especificar'
$arr[91,8] = 'Primera'
$arr[91,9] = 1000
$arr[91,10] = 700
$arr[91,11] = 800
$arr[91,12] = 750
$arr[91,13] = '$/unidad'
$arr[91,14] = 'Región del Maule'
$arr[91,15] = 750
$arr[91,16] = 1
$arr[91,17] = 'Hortaliza'
$arr[92,0] = 11
$arr[92,1] = 'Vega Monumental Concepción'
$arr[92,2] = 'Bíobío'
$arr[92,3] = 44196
$arr[92,4] = 8
$arr[92,5] = 100112023
$arr[92,6] = 'Brócoli'
$arr[92,7] = 'Sin especificar'
$arr[92,8] = 'Segunda'
$arr[92,9] = 500
$arr[92,10] = 600
$arr[92,11] = 600
$arr[92,12] = 600
$arr[92,13] = '$/unidad'
$arr[92,14] = 'Región del Maule'
$arr[92,15] = 600
$arr[92,16] = 1
$arr[92,17] = 'Hortaliza'
$arr[93,0] = 11
$arr[93,1] = 'Vega Monumental Concepción'
$arr[93,2] = 'Bíobío'
$arr[93,3] = 44496
$arr[93,4] = 8
$arr[93,5] = 100112023
$arr[93,6] = 'Brócoli'
$arr[93,7] = 'Sin especificar'
$arr[93,8] = 'Primera'
$arr[93,9] = 1050
$arr[93,10] = 750
$arr[93,11] = 750
$arr[93,12] = 750
$arr[93,13] = '$/unidad'
$arr[93,14] = 'Región del Maule'
$arr[93,15] = 750
$arr[93,16] = 1
$arr[93,17] = 'Hortaliza'
$arr[94,0] = 11
$arr[94,1] = 'Vega Monumental Concepción'
$arr[94,2] = 'Bíobío'
$arr[94,3] = 44425
$arr[94,4] = 8
$arr[94,5] = 100112023
$arr[94,6] = 'Brócoli'
$arr[94,7] = 'Sin especificar'
$arr[94,8] = 'Primera'
$arr[94,9] = 2000
$arr[94,10] = 600
$arr[94,11] = 700
$arr[94,12] = 650
$arr[94,13] = '$/unidad'
$arr[94,14] = 'Región Metropolitana'
$arr[94,15] = 650
$arr[94,16] = 1
$arr[94,17] = 'Hortaliza'
$arr[95,0] = 11
$arr[95,1] = 'Vega Monumental Concepción'
$arr[95,2] = 'Bíobío'
$arr[95,3] = 44425
$arr[95,4] = 8
$arr[95,5] = 100112023
$arr[95,6] = 'Brócoli'
$arr[95,7] = 'Sin especificar'
$arr[95,8] = 'Segunda'
$arr[95,9] = 1000
$arr[95,10] = 500
$arr[95,11] = 500
$arr[95,12] = 500
$arr[95,13] = '$/unidad'
$arr[95,14] = 'Región Metropolitana'
$arr[95,15] = 500
$arr[95,16] = 1
$arr[95,17] = 'Hortaliza'
$arr[96,0] = 11
$arr[96,1] = 'Vega Monumental Concepción'
$arr[96,2] = 'Bíobío'
$arr[96,3] = 44512
$arr[96,4] = 8
$arr[96,5] = 100112023
$arr[96,6] = 'Brócoli'
$arr[96,7] = 'Sin especificar'
$arr[96,8] = 'Primera'
$arr[96,9] = 2000
$arr[96,10] = 700
$arr[96,11] = 800
$arr[96,12] = 750
$arr[96,13] = '$/unidad'
$arr[96,14] = 'Región Metropolitana'
$arr[96,15] = 750
$arr[96,16] = 1
$arr[96,17] = 'Hortaliza'
$arr[97,0] = 11
$arr[97,1] = 'Vega Monumental Concepción'
$arr[97,2] = 'Bíobío'
$arr[97,3] = 44512
$arr[97,4] = 8
$arr[97,5] = 100112023
$arr[97,6] = 'Brócoli'
$arr[97,7] = 'Sin especificar'
$arr[97,8] = 'Segunda'
$arr[97,9] = 1000
$arr[97,10] = 600
$arr[97,11] = 600
$arr[97,12] = 600
$arr[97,13] = '$/unidad'
$arr[97,14] = 'Región Metropolitana'
$arr[97,15] = 600
$arr[97,16] = 1
$arr[97,17] = 'Hortaliza'
$arr[98,0] = 11
$arr[98,1] = 'Vega Monumental Concepción'
$arr[98,2] = 'Bíobío'
$arr[98,3] = 44497
$arr[98,4] = 8
$arr[98,5] = 100112023
$arr[98,6] = 'Brócoli'
$arr[98,7] = 'Sin especificar'
$arr[98,8] = 'Primera'
$arr[98,9] = 2000
$arr[98,10] = 650
$arr[98,11] = 700
$arr[98,12] = 675
$arr[98,13] = '$/unidad'
$arr[98,14] = 'Región Metropolitana'
$arr[98,15] = 675
$arr[98,16] = 1
$arr[98,17] = 'Hortaliza'
$arr[99,0] = 11
$arr[99,1] = 'Vega Monumental Concepción'
$arr[99,2] = 'Bíobío'
$arr[99,3] = 44497
$arr[99,4] = 8
$arr[99,5] = 100112023
$arr[99,6] = 'Brócoli'
$arr[99,7] = 'Sin especificar'
$arr[99,8] = 'Segunda'
$arr[99,9] = 1500
$arr[99,10] = 600
$arr[99,11] = 600
$arr[99,12] = 600
$arr[99,13] = '$/unidad'
$arr[99,14] = 'Región Metropolitana'
$arr[99,15] = 600
$arr[99,16] = 1
$arr[99,17] = 'Hortaliza'
$arr[100,0] = 11
$arr[100,1] = 'Vega Monumental Concepción'
$arr[100,2] = 'Bíobío'
$arr[100,3] = 44362
$arr[100,4] = 8
$arr[100,5] = 100112023
$arr[100,6] = 'Brócoli'
$arr[100,7] = 'Sin especificar'
$arr[100,8] = 'Primera'
$arr[100,9] = 2000
$arr[100,10] = 600
$arr[100,11] = 700
$arr[100,12] = 650
$arr[100,13] = '$/unidad'
$arr[100,14] = 'Región Metropolitana'
$arr[100,15] = 650
$arr[100,16] = 1
$arr[100,17] = 'Hortaliza'
$arr[101,0] = 11
$arr[101,1] = 'Vega Monumental Concepción'
$arr[101,2] = 'Bíobío'
$arr[101,3] = 44362
$arr[101,4] = 8
$arr[101,5] = 100112023
$arr[101,6] = 'Brócoli'
$arr[101,7] = 'Sin especificar'
$arr[101,8] = 'Segunda'
$arr[101,9] = 1000
$arr[101,10] = 500
$arr[101,11] = 500
$arr[101,12] = 500
$arr[101,13] = '$/unidad'
$arr[101,14] = 'Región Metropolitana'
$arr[101,15] = 500
$arr[101,16] = 1
$arr[101,17] = 'Hortaliza'
$arr[102,0] = 11
$arr[102,1] = 'Vega Monumental Concepción'
$arr[102,2] = 'Bíobío'
$arr[102,3] = 44747
$arr[102,4] = 8
$arr[102,5] = 100112023
$arr[102,6] = 'Brócoli'
$arr[102,7] = 'Sin especificar'
$arr[102,8] = 'Primera'
$arr[102,9] = 2000
$arr[102,10] = 800
$arr[102,11] = 900
$arr[102,12] = 850
$arr[102,13] = '$/unidad'
$arr[102,14] = 'Región Metropolitana'
$arr[102,15] = 850
$arr[102,16] = 1
$arr[102,17] = 'Hortaliza'
$arr[103,0] = 11
$arr[103,1] = 'Vega Monumental Concepción'
$arr[103,2] = 'Bíobío'
$arr[103,3] = 44357
$arr[103,4] = 8
$arr[103,5] = 100112023
$arr[103,6] = 'Brócoli'
$arr[103,7] = 'Sin especificar'
$arr[103,8] = 'Primera'
$arr[103,9] = 2000
$arr[103,10] = 600
$arr[103,11] = 700
$arr[103,12] = 650
$arr[103,13] = '$/unidad'
$arr[103,14] = 'Región Metropolitana'
$arr[103,15] = 650
$arr[103,16] = 1
$arr[103,17] = 'Hortaliza'
$arr[104,0] = 11
$arr[104,1] = 'Vega Monumental Concepción'
$arr[104,2] = 'Bíobío'
$arr[104,3] = 44357
$arr[104,4] = 8
$arr[104,5] = 100112023
$arr[104,6] = 'Brócoli'
$arr[104,7] = 'Sin especificar'
$arr[104,8] = 'Segunda'
$arr[104,9] = 1000
$arr[104,10] = 500
$arr[104,11] = 500
$arr[104,12] = 500
$arr[104,13] = '$/unidad'
$arr[104,14] = 'Región Metropolitana'
$arr[104,15] = 500
$arr[104,16] = 1
$arr[104,17] = 'Hortaliza'
$arr[105,0] = 11
$arr[105,1] = 'Vega Monumental Concepción'
$arr[105,2] = 'Bíobío'
$arr[105,3] = 44279
$arr[105,4] = 8
$arr[105,5] = 100112023
$arr[105,6] = 'Brócoli'
$arr[105,7] = 'Sin especificar'
$arr[105,8] = 'Primera'
$arr[105,9] = 1000
$arr[105,10] = 700
$arr[105,11] = 800
$arr[105,12] = 750
$arr[105,13] = '$/unidad'
$arr[105,14] = 'Región del Maule'
$arr[105,15] = 750
$arr[105,16] = 1
$arr[105,17] = 'Hortaliza'
$arr[106,0] = 11
$arr[106,1] = 'Vega Monumental Concepción'
$arr[106,2] = 'Bíobío'
$arr[106,3] = 44279
$arr[106,4] = 8
$arr[106,5] = 100112023
$arr[106,6] = 'Brócoli'
$arr[106,7] = 'Sin especificar'
$arr[106,8] = 'Segunda'
$arr[106,9] = 500
$arr[106,10] = 600
$arr[106,11] = 600
$arr[106,12] = 600
$arr[106,13] = '$/unidad'
$arr[106,14] = 'Región del Maule'
$arr[106,15] = 600
$arr[106,16] = 1
$arr[106,17] = 'Hortaliza'
$arr[107,0] = 11
$arr[107,1] = 'Vega Monumental Concepción'
$arr[107,2] = 'Bíobío'
$arr[107,3] = 44551
$arr[107,4] = 8
$arr[107,5] = 100112023
$arr[107,6] = 'Brócoli'
$arr[107,7] = 'Sin especificar'
$arr[107,8] = 'Primera'
$arr[107,9] = 2000
$arr[107,10] = 600
$arr[107,11] = 700
$arr[107,12] = 650
$arr[107,13] = '$/unidad'
$arr[107,14] = 'Región Metropolitana'
$arr[107,15] = 650
$arr[107,16] = 1
$arr[107,17] = 'Hortaliza'
$arr[108,0] = 11
$arr[108,1] = 'Vega Monumental Concepción'
$arr[108,2] = 'Bíobío'
$arr[108,3] = 44551
$arr[108,4] = 8
$arr[108,5] = 100112023
$arr[108,6] = 'Brócoli'
$arr[108,7] = 'Sin especificar'
$arr[108,8] = 'Segunda'
$arr[108,9] = 1000
$arr[108,10] = 500
$arr[108,11] = 500
$arr[108,12] = 500
$arr[108,13] = '$/unidad'
$arr[108,14] = 'Región Metropolitana'
$arr[108,15] = 500
$arr[108,16] = 1
$arr[108,17] = 'Hortaliza'
$arr[109,0] = 11
$arr[109,1] = 'Vega Monumental Concepción'
$arr[109,2] = 'Bíobío'
$arr[109,3] = 44547
$arr[109,4] = 8
$arr[109,5] = 100112023
$arr[109,6] = 'Brócoli'
$arr[109,7] = 'Sin especificar'
$arr[109,8] = 'Primera'
$arr[109,9] = 1000
$arr[109,10] = 600
$arr[109,11] = 700
$arr[109,12] = 650
$arr[109,13] = '$/unidad'
$arr[109,14] = 'Región Metropolitana'
$arr[109,15] = 650
$arr[109,16] = 1
$arr[109,17] = 'Hortaliza'
$arr[110,0] = 11
$arr[110,1] = 'Vega Monumental Concepción'
$arr[110,2] = 'Bíobío'
$arr[110,3] = 44547
$arr[110,4] = 8
$arr[110,5] = 100112023
$arr[110,6] = 'Brócoli'
$arr[110,7] = 'Sin especificar'
$arr[110,8] = 'Segunda'
$arr[110,9] = 500
$arr[110,10] = 500
$arr[110,11] = 500
$arr[110,12] = 500
$arr[110,13] = '$/unidad'
$arr[110,14] = 'Región Metropolitana'
$arr[110,15] = 500
$arr[110,16] = 1
$arr[110,17] = 'Hortaliza'
$arr[111,0] = 11
$arr[111,1] = 'Vega Monumental Concepción'
$arr[111,2] = 'Bíobío'
$arr[111,3] = 44355
$arr[111,4] = 8
$arr[111,5] = 100112023
$arr[111,6] = 'Brócoli'
$arr[111,7] = 'Sin especificar'
$arr[111,8] = 'Primera'
$arr[111,9] = 2000
$arr[111,10] = 600
$arr[111,11] = 700
$arr[111,12] = 650
$arr[111,13] = '$/unidad'
$arr[111,14] = 'Región Metropolitana'
$arr[111,15] = 650
$arr[111,16] = 1
$arr[111,17] = 'Hortaliza'
$arr[112,0] = 11
$arr[112,1] = 'Vega Monumental Concepción'
$arr[112,2] = 'Bíobío'
$arr[112,3] = 44355
$arr[112,4] = 8
$arr[112,5] = 100112023
$arr[112,6] = 'Brócoli'
$arr[112,7] = 'Sin especificar'
$arr[112,8] = 'Segunda'
$arr[112,9] = 1000
$arr[112,10] = 500
$arr[112,11] = 500
$arr[112,12] = 500
$arr[112,13] = '$/unidad'
$arr[112,14] = 'Región Metropolitana'
$arr[112,15] = 500
$arr[112,16] = 1
$arr[112,17] = 'Hortaliza'
$arr[113,0] = 11
$arr[113,1] = 'Vega Monumental Concepción'
$arr[113,2] = 'Bíobío'
$arr[113,3] = 44657
$arr[113,4] = 8
$arr[113,5] = 100112023
$arr[113,6] = 'Brócoli'
$arr[113,7] = 'Sin especificar'
$arr[113,8] = 'Primera'
$arr[113,9] = 1000
$arr[113,10] = 800
$arr[113,11] = 1000
$arr[113,12] = 900
$arr[113,13] = '$/unidad'
$arr[113,14] = 'Región Metropolitana'
$arr[113,15] = 900
$arr[113,16] = 1
$arr[113,17] = 'Hortaliza'
$arr[114,0] = 11
$arr[114,1] = 'Vega Monumental Concepción'
$arr[114,2] = 'Bíobío'
$arr[114,3] = 44391
$arr[114,4] = 8
$arr[114,5] = 100112023
$arr[114,6] = 'Brócoli'
$arr[114,7] = 'Sin especificar'
$arr[114,8] = 'Primera'
$arr[114,9] = 2000
$arr[114,10] = 600
$arr[114,11] = 700
$arr[114,12] = 650
$arr[114,13] = '$/unidad'
$arr[114,14] = 'Región Metropolitana'
$arr[114,15] = 650
$arr[114,16] = 1
$arr[114,17] = 'Hortaliza'
$arr[115,0] = 11
$arr[115,1] = 'Vega Monumental Concepción'
$arr[115,2] = 'Bíobío'
$arr[115,3] = 44391
$arr[115,4] = 8
$arr[115,5] = 100112023
$arr[115,6] = 'Brócoli'
$arr[115,7] = 'Sin especificar'
$arr[115,8] = 'Segunda'
$arr[115,9] = 1000
$arr[115,10] = 500
$arr[115,11] = 500
$arr[115,12] = 500
$arr[115,13] = '$/unidad'
$arr[115,14] = 'Región Metropolitana'
$arr[115,15] = 500
$arr[115,16] = 1
$arr[115,17] = 'Hortaliza'
$arr[116,0] = 11
$arr[116,1] = 'Vega Monumental Concepción'
$arr[116,2] = 'Bíobío'
$arr[116,3] = 44453
$arr[116,4] = 8
$arr[116,5] = 100112023
$arr[116,6] = 'Brócoli'
$arr[116,7] = 'Sin especificar'
$arr[116,8] = 'Primera'
$arr[116,9] = 2000
$arr[116,10] = 600
$arr[116,11] = 700
$arr[116,12] = 650
$arr[116,13] = '$/unidad'
$arr[116,14] = 'Región Metropolitana'
$arr[116,15] = 650
$arr[116,16] = 1
$arr[116,17] = 'Hortaliza'
$arr[117,0] = 11
$arr[117,1] = 'Vega Monumental Concepción'
$arr[117,2] = 'Bíobío'
$arr[117,3] = 44453
$arr[117,4] = 8
$arr[117,5] = 100112023
$arr[117,6] = 'Brócoli'
$arr[117,7] = 'Sin especificar'
$arr[117,8] = 'Segunda'
$arr[117,9] = 1000
$arr[117,10] = 500
$arr[117,11] = 500
$arr[117,12] = 500
$arr[117,13] = '$/unidad'
$arr[117,14] = 'Región Metropolitana'
$arr[117,15] = 500
$arr[117,16] = 1
$arr[117,17] = 'Hortaliza'
$arr[118,0] = 11
$arr[118,1] = 'Vega Monumental Concepción'
$arr[118,2] = 'Bíobío'
$arr[118,3] = 44189
$arr[118,4] = 8
$arr[118,5] = 100112023
$arr[118,6] = 'Brócoli'
$arr[118,7] = 'Sin especificar'
$arr[118,8] = 'Primera'
$arr[118,9] = 800
$arr[118,10] = 600
$arr[118,11] = 700
$arr[118,12] = 650
$arr[118,13] = '$/unidad'
$arr[118,14] = 'Región Metropolitana'
$arr[118,15] = 650
$arr[118,16] = 1
$arr[118,17] = 'Hortaliza'
$arr[119,0] = 11
$arr[119,1] = 'Vega Monumental Concepción'
$arr[119,2] = 'Bíobío'
$arr[119,3] = 44189
$arr[119,4] = 8
$arr[119,5] = 100112023
$arr[119,6] = 'Brócoli'
$arr[119,7] = 'Sin especificar'
$arr[119,8] = 'Segunda'
$arr[119,9] = 400
$arr[119,10] = 500
$arr[119,11] = 500
$arr[119,12] = 500
$arr[119,13] = '$/unidad'
$arr[119,14] = 'Región Metropolitana'
$arr[119,15] = 500
$arr[119,16] = 1
$arr[119,17] = 'Hortaliza'
$arr[120,0] = 11
$arr[120,1] = 'Vega Monumental Concepción'
$arr[120,2] = 'Bíobío'
$arr[120,3] = 44358
$arr[120,4] = 8
$arr[120,5] = 100112023
$arr[120,6] = 'Brócoli'
$arr[120,7] = 'Sin especificar'
$arr[120,8] = 'Primera'
$arr[120,9] = 2000
$arr[120,10] = 500
$arr[120,11] = 600
$arr[120,12] = 550
$arr[120,13] = '$/unidad'
$arr[120,14] = 'Región del Maule'
$arr[120,15] = 550
$arr[120,16] = 1
$arr[120,17] = 'Hortaliza'
$arr[121,0] = 11
$arr[121,1] = 'Vega Monumental Concepción'
$arr[121,2] = 'Bíobío'
$arr[121,3] = 44358
$arr[121,4] = 8
$arr[121,5] = 100112023
$arr[121,6] = 'Brócoli'
$arr[121,7] = 'Sin especificar'
$arr[121,8] = 'Segunda'
$arr[121,9] = 1000
$arr[121,10] = 400
$arr[121,11] = 400
$arr[121,12] = 400
$arr[121,13] = '$/unidad'
$arr[121,14] = 'Región del Maule'
$arr[121,15] = 400
$arr[121,16] = 1
$arr[121,17] = 'Hortaliza'
$arr[122,0] = 11
$arr[122,1] = 'Vega Monumental Concepción'
$arr[122,2] = 'Bíobío'
$arr[122,3] = 44572
$arr[122,4] = 8
$arr[122,5] = 100112023
$arr[122,6] = 'Brócoli'
$arr[122,7] = 'Sin especificar'
$arr[122,8] = 'Primera'
$arr[122,9] = 270
$arr[122,10] = 850
$arr[122,11] = 900
$arr[122,12] = 878
$arr[122,13] = '$/unidad'
$arr[122,14] = 'Región Metropolitana'
$arr[122,15] = 878
$arr[122,16] = 1
$arr[122,17] = 'Hortaliza'

# Make sure the two brand-new rows at the bottom (326-327) pick up
# the same date number-format used by the rest of column D before
# writing the values (new cells otherwise default to General).
$ws.Range("D326:D327").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A205:R327").Value = $arr

Write-Output "applied rows 205-327"
